$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three more studies left to analyze -> add rows 7-11 (Tag + "Complete?" = Yes, in the
# order they were typed in: Schkade2010, Myers1975, Moscovici1969, Krizan2007, Myers1970).
$ws.Range("A11").Value = "Schkade2010"
$ws.Range("B11").Value = "Yes"

$ws.Range("A10").Value = "Myers1975"
$ws.Range("B10").Value = "Yes"

$ws.Range("A8").Value = "Moscovici1969"
$ws.Range("B8").Value = "Yes"

$ws.Range("A7").Value = "Krizan2007"

$ws.Range("A9").Value = "Myers1970"
$ws.Range("B9").Value = "Yes"

# The conditional-formatting ranges for columns B and C were anchored through row 12;
# now that the data only runs through row 11, pull the ranges back in by one row.
$cfs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $rule = $cfs.Item($i)
    $addr = $rule.AppliesTo.Address()
    if ($addr -eq '$B$3:$B$12') {
        $rule.ModifyAppliesToRange($ws.Range("B3:B11"))
    } elseif ($addr -eq '$C$2:$C$12') {
        $rule.ModifyAppliesToRange($ws.Range("C2:C11"))
    }
}

$ws.Range("B7").Select()
